$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "42.947.24"
$ws.Cells.Item(2, 5).Value = "  +4.08%  "

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "2.228.33"
$ws.Cells.Item(3, 5).Value = "  +3.49%  "

$ws.Cells.Item(4, 5).Value = "  -0.04%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "252.11"
$ws.Cells.Item(5, 5).Value = "  +6.60%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.616"
$ws.Cells.Item(6, 5).Value = "  +2.28%  "

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "75.06"
$ws.Cells.Item(7, 5).Value = "  +6.78%  "

$ws.Cells.Item(8, 5).Value = "  +0.00%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.595"
$ws.Cells.Item(9, 5).Value = "  +4.22%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "41.09"
$ws.Cells.Item(10, 5).Value = "  +4.22%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0920"
$ws.Cells.Item(11, 5).Value = "  +2.58%  "

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "6.86"
$ws.Cells.Item(12, 5).Value = "  +3.53%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.102"
$ws.Cells.Item(13, 5).Value = "  +1.83%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "2.564.38"
$ws.Cells.Item(14, 5).Value = "  +3.58%  "

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "14.49"
$ws.Cells.Item(15, 5).Value = "  +2.52%  "

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "2.230.49"
$ws.Cells.Item(16, 5).Value = "  +4.49%  "

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "0.788"
$ws.Cells.Item(17, 5).Value = "  +0.90%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "42.873.41"
$ws.Cells.Item(18, 5).Value = "  +4.08%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.0000103"
$ws.Cells.Item(19, 5).Value = "  +3.27%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "71.19"
$ws.Cells.Item(20, 5).Value = "  +2.92%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "5.93"
$ws.Cells.Item(21, 5).Value = "  +3.59%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "229.92"
$ws.Cells.Item(22, 5).Value = "  +2.09%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "2.18"
$ws.Cells.Item(23, 5).Value = "  +10.57%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "9.29"
$ws.Cells.Item(24, 5).Value = "  -3.14%  "

$ws.Cells.Item(25, 5).Value = "  +0.08%  "

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "10.68"
$ws.Cells.Item(26, 5).Value = "  +0.76%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "3.41"
$ws.Cells.Item(27, 5).Value = "  +4.41%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "39.07"
$ws.Cells.Item(28, 5).Value = "  +22.90%  "

$ws.Cells.Item(29, 2).Value = "Toncoin"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "2.23"
$ws.Cells.Item(29, 5).Value = "  +3.05%  "

$ws.Cells.Item(30, 2).Value = "PancakeSwap"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "2.21"
$ws.Cells.Item(30, 5).Value = "  +2.24%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "170.20"
$ws.Cells.Item(31, 5).Value = "  -0.45%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "20.15"
$ws.Cells.Item(32, 5).Value = "  +2.48%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.0793"
$ws.Cells.Item(33, 5).Value = "  +4.52%  "

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "5.23"
$ws.Cells.Item(34, 5).Value = "  +3.17%  "

$ws.Cells.Item(35, 5).Value = "  +9.43%  "

$ws.Cells.Item(36, 5).Value = "  +0.68%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "4.45"
$ws.Cells.Item(37, 5).Value = "  +3.49%  "

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.0324"
$ws.Cells.Item(38, 5).Value = "  +11.54%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "12.20"
$ws.Cells.Item(39, 5).Value = "  +1.61%  "

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "2.10"
$ws.Cells.Item(40, 5).Value = "  +2.53%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.204"
$ws.Cells.Item(41, 5).Value = "  +8.52%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "5.35"
$ws.Cells.Item(42, 5).Value = "  +2.18%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "59.38"
$ws.Cells.Item(43, 5).Value = "  +1.99%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "8.62"
$ws.Cells.Item(44, 5).Value = "  +4.59%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "103.16"
$ws.Cells.Item(45, 5).Value = "  +6.02%  "

$ws.Cells.Item(46, 2).Value = "Cronos"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.0984"
$ws.Cells.Item(46, 5).Value = "  +2.30%  "

$ws.Cells.Item(47, 2).Value = "WOONetwork"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.469"
$ws.Cells.Item(47, 5).Value = "  +20.92%  "

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "2.42"
$ws.Cells.Item(48, 5).Value = "  +12.48%  "

$ws.Cells.Item(49, 5).Value = "  +3.81%  "

$ws.Cells.Item(50, 5).Value = "  +2.88%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "2.68"
$ws.Cells.Item(51, 5).Value = "  +2.40%  "
